$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data: plate/reagent info stays constant; each reagent (source well D / reagent I / volume H)
# is now dispensed into both destination wells A1 and A2 (column G), doubling rows 2-4 into rows 2-7.

$sourcePlateName = "level 2 LDV source plate"
$sourcePlateType = "384LDV_AQ_B"
$destPlateName   = "384-Well Level 2 MoClo output plate"
$destPlateType   = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"

$reagents = @(
    @{ SourceWell = "A1"; Volume = 500; Reagent = "DNA ligase buffer" },
    @{ SourceWell = "A2"; Volume = 125; Reagent = "DNA ligase" },
    @{ SourceWell = "A3"; Volume = 250; Reagent = "BsmBI (NEB)" }
)
$destWells = @("A1", "A2")

$uid = 1
$row = 2
foreach ($reagent in $reagents) {
    foreach ($destWell in $destWells) {
        $ws.Cells.Item($row, 1).Value = $uid
        $ws.Cells.Item($row, 2).Value = $sourcePlateName
        $ws.Cells.Item($row, 3).Value = $sourcePlateType
        $ws.Cells.Item($row, 4).Value = $reagent.SourceWell
        $ws.Cells.Item($row, 5).Value = $destPlateName
        $ws.Cells.Item($row, 6).Value = $destPlateType
        $ws.Cells.Item($row, 7).Value = $destWell
        $ws.Cells.Item($row, 8).Value = $reagent.Volume
        $ws.Cells.Item($row, 9).Value = $reagent.Reagent

        $uid = $uid + 1
        $row = $row + 1
    }
}
